$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the first sheet: 每日一题 -> 题库
#    (the _xlnm._FilterDatabase defined name follows the sheet automatically)
# ---------------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item(1)
$sheet1.Name = "题库"
$sheet2 = $wb.Worksheets.Item(2)
$sheet3 = $wb.Worksheets.Item(3)
$sheet4 = $wb.Worksheets.Item(4)

# ---------------------------------------------------------------------------
# 2. 题库 sheet (sheet1.xml): tag LC 2262 with a new note, then append the
#    5 new rows (795-799) for Weekly Contest 329.
# ---------------------------------------------------------------------------

# G365 gets a new wrap-text note (style like G2 / G794, font 6 + wrapText)
$sheet1.Range("G2").Copy()
$sheet1.Range("G365").PasteSpecial(-4122)
$sheet1.Range("G365").Value = "[DP]一维DP，dp[i]表示以s[i]结尾的某种性质的计数"

# Row 795: 交替数字和 (简单)
$sheet1.Range("B794").Copy()
$sheet1.Range("B795").PasteSpecial(-4122)
$sheet1.Range("B795").Value = "交替数字和"

$sheet1.Range("C794").Copy()
$sheet1.Range("C795").PasteSpecial(-4122)
$sheet1.Range("C795").Value = 44948

$sheet1.Range("D794").Copy()
$sheet1.Range("D795").PasteSpecial(-4122)
$sheet1.Range("D795").Value = "简单"

# Row 796: 根据第 K 场考试的分数排序 (中等)
$sheet1.Range("B794").Copy()
$sheet1.Range("B796").PasteSpecial(-4122)
$sheet1.Range("B796").Value = "根据第 K 场考试的分数排序"

$sheet1.Range("C794").Copy()
$sheet1.Range("C796").PasteSpecial(-4122)
$sheet1.Range("C796").Value = 44948

$sheet1.Range("D794").Copy()
$sheet1.Range("D796").PasteSpecial(-4122)
$sheet1.Range("D796").Value = "中等"

$sheet1.Range("G2").Copy()
$sheet1.Range("F796").PasteSpecial(-4122)
$sheet1.Range("F796").Value = "排序"

$sheet1.Range("G2").Copy()
$sheet1.Range("G796").PasteSpecial(-4122)
$sheet1.Range("G796").Value = "用enumrate的排序"

# Row 797: 执行逐位运算使字符串相等 (中等)
$sheet1.Range("B794").Copy()
$sheet1.Range("B797").PasteSpecial(-4122)
$sheet1.Range("B797").Value = "执行逐位运算使字符串相等"

$sheet1.Range("C794").Copy()
$sheet1.Range("C797").PasteSpecial(-4122)
$sheet1.Range("C797").Value = 44948

$sheet1.Range("D794").Copy()
$sheet1.Range("D797").PasteSpecial(-4122)
$sheet1.Range("D797").Value = "中等"

$sheet1.Range("G2").Copy()
$sheet1.Range("F797").PasteSpecial(-4122)
$sheet1.Range("F797").Value = "数学"

# Row 798: 拆分数组的最小代价 (困难)
$sheet1.Range("B794").Copy()
$sheet1.Range("B798").PasteSpecial(-4122)
$sheet1.Range("B798").Value = "拆分数组的最小代价"

$sheet1.Range("C794").Copy()
$sheet1.Range("C798").PasteSpecial(-4122)
$sheet1.Range("C798").Value = 44948

$sheet1.Range("D794").Copy()
$sheet1.Range("D798").PasteSpecial(-4122)
$sheet1.Range("D798").Value = "困难"

$sheet1.Range("G2").Copy()
$sheet1.Range("G798").PasteSpecial(-4122)
$sheet1.Range("G798").Value = "[子数组]"

# Row 799: 统计特殊子序列的数目 (困难), LC 1955
$sheet1.Range("A794").Copy()
$sheet1.Range("A799").PasteSpecial(-4122)
$sheet1.Range("A799").Value = 1955

$sheet1.Range("B794").Copy()
$sheet1.Range("B799").PasteSpecial(-4122)
$sheet1.Range("B799").Value = "统计特殊子序列的数目"

$sheet1.Range("C794").Copy()
$sheet1.Range("C799").PasteSpecial(-4122)
$sheet1.Range("C799").Value = 44948

$sheet1.Range("D794").Copy()
$sheet1.Range("D799").PasteSpecial(-4122)
$sheet1.Range("D799").Value = "困难"

$sheet1.Range("E794").Copy()
$sheet1.Range("E799").PasteSpecial(-4122)
$sheet1.Range("E799").Value = 2125

$sheet1.Range("G2").Copy()
$sheet1.Range("F799").PasteSpecial(-4122)
$sheet1.Range("F799").Value = "DP"

$sheet1.Range("G2").Copy()
$sheet1.Range("G799").PasteSpecial(-4122)
$sheet1.Range("G799").Value = "[小二维DP]，实现是可以压缩掉一维"

# Keep the frozen header pane, scroll the view to show the new rows.
$sheet1.Activate()
$sheet1.Range("A778").Select()
$excel.ActiveWindow.FreezePanes = $true
$sheet1.Range("G799").Select()

# ---------------------------------------------------------------------------
# 3. 算法分类 sheet (sheet2.xml): add two rows under the 树形DP section,
#    column H.
# ---------------------------------------------------------------------------
$sheet2.Range("C14").Copy()
$sheet2.Range("H13").PasteSpecial(-4122)
$sheet2.Range("H13").Value = "树形DP"

$sheet2.Range("C14").Copy()
$sheet2.Range("H14").PasteSpecial(-4122)
$sheet2.Range("H14").Value = "0-1 BFS"

$sheet2.Activate()
$sheet2.Range("H15").Select()

# ---------------------------------------------------------------------------
# 4. Sheet3 (Sheet3): insert a new row 10 (小二维DP) before the existing
#    双指针 row, shifting everything below down by one.
#
# Original layout:
#   row10: D=双指针
#   row11: D=一维DP
#   row12: D=单调栈+二分, E=300,  F=LIS问题可以用这个来优化
#   row13: D=前后缀统计,  E=2484, F=回文子序列
#   row14: D=DP+线段树优化, E=2407, F=LIS
#
# Target layout:
#   row10 (NEW): D=小二维DP, E=1955, F=只需要保存dp[i]和dp[i-1]两项即可, H/J blank-styled
#   row11: D=双指针 (ht 28.8), H=DP, I=2262, J=[DP]一维DP...note
#   row12: D=一维DP
#   row13: D=单调栈+二分, E=300,  F=LIS问题可以用这个来优化
#   row14: D=前后缀统计,  E=2484, F=回文子序列
#   row15: D=DP+线段树优化, E=2407, F=LIS
# ---------------------------------------------------------------------------

# Snapshot the old values BEFORE writing anything, so shifting rows down
# never reads back data we already overwrote.
$old_D10 = $sheet3.Range("D10").Value()
$old_D11 = $sheet3.Range("D11").Value()
$old_D12 = $sheet3.Range("D12").Value()
$old_E12 = $sheet3.Range("E12").Value()
$old_F12 = $sheet3.Range("F12").Value()
$old_D13 = $sheet3.Range("D13").Value()
$old_E13 = $sheet3.Range("E13").Value()
$old_F13 = $sheet3.Range("F13").Value()
$old_D14 = $sheet3.Range("D14").Value()
$old_E14 = $sheet3.Range("E14").Value()
$old_F14 = $sheet3.Range("F14").Value()

# row15 (was row14): D=DP+线段树优化, E=2407, F=LIS
$sheet3.Range("D14").Copy()
$sheet3.Range("D15").PasteSpecial(-4122)
$sheet3.Range("D15").Value = $old_D14
$sheet3.Range("E14").Copy()
$sheet3.Range("E15").PasteSpecial(-4122)
$sheet3.Range("E15").Value = $old_E14
$sheet3.Range("F14").Copy()
$sheet3.Range("F15").PasteSpecial(-4122)
$sheet3.Range("F15").Value = $old_F14

# row14 (was row13): D=前后缀统计, E=2484, F=回文子序列
$sheet3.Range("D13").Copy()
$sheet3.Range("D14").PasteSpecial(-4122)
$sheet3.Range("D14").Value = $old_D13
$sheet3.Range("E13").Copy()
$sheet3.Range("E14").PasteSpecial(-4122)
$sheet3.Range("E14").Value = $old_E13
$sheet3.Range("F13").Copy()
$sheet3.Range("F14").PasteSpecial(-4122)
$sheet3.Range("F14").Value = $old_F13

# row13 (was row12): D=单调栈+二分, E=300, F=LIS问题可以用这个来优化
$sheet3.Range("D12").Copy()
$sheet3.Range("D13").PasteSpecial(-4122)
$sheet3.Range("D13").Value = $old_D12
$sheet3.Range("E12").Copy()
$sheet3.Range("E13").PasteSpecial(-4122)
$sheet3.Range("E13").Value = $old_E12
$sheet3.Range("F12").Copy()
$sheet3.Range("F13").PasteSpecial(-4122)
$sheet3.Range("F13").Value = $old_F12

# row12 (was row11): D=一维DP
$sheet3.Range("D11").Copy()
$sheet3.Range("D12").PasteSpecial(-4122)
$sheet3.Range("D12").Value = $old_D11

# row11 (was row10): D=双指针, plus new DP example note
$sheet3.Range("D10").Copy()
$sheet3.Range("D11").PasteSpecial(-4122)
$sheet3.Range("D11").Value = $old_D10

$sheet3.Range("H9").Copy()
$sheet3.Range("H11").PasteSpecial(-4122)
$sheet3.Range("H11").Value = "DP"
$sheet3.Range("I9").Copy()
$sheet3.Range("I11").PasteSpecial(-4122)
$sheet3.Range("I11").Value = 2262
$sheet1.Range("G2").Copy()
$sheet3.Range("J11").PasteSpecial(-4122)
$sheet3.Range("J11").Value = "[DP]一维DP，dp[i]表示以s[i]结尾的某种性质的计数" + [char]10 + "最终结果为 sum(dp)"
$sheet3.Rows.Item(11).RowHeight = 28.8

# row10 (brand new): 小二维DP example
$sheet3.Range("D10").Value = "小二维DP"
$sheet3.Range("E9").Copy()
$sheet3.Range("E10").PasteSpecial(-4122)
$sheet3.Range("E10").Value = 1955
$sheet3.Range("D9").Copy()
$sheet3.Range("F10").PasteSpecial(-4122)
$sheet3.Range("F10").Value = "只需要保存dp[i] 和 dp[i-1] 两项即可"
$sheet3.Range("H9").Copy()
$sheet3.Range("H10").PasteSpecial(-4122)
$sheet3.Range("H10").ClearContents()
$sheet3.Range("J9").Copy()
$sheet3.Range("J10").PasteSpecial(-4122)
$sheet3.Range("J10").ClearContents()

$sheet3.Columns.Item(10).AutoFit()

$sheet3.Activate()
$sheet3.Range("F21").Select()

# ---------------------------------------------------------------------------
# 5. 题型分类 sheet (sheet4.xml): add a new row 11 for the sorted() tip.
# ---------------------------------------------------------------------------
$sheet4.Range("D9").Copy()
$sheet4.Range("E11").PasteSpecial(-4122)
$sheet4.Range("E11").Value = "sorted会拷贝一份数组，空间复杂度比sort要大"

$sheet4.Activate()
$sheet4.Range("E11").Select()

$sheet1.Activate()
